# Fix circular reference in Cash Flow Statement
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cash Flow Statement")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")

foreach ($col in $cols) {
    # Row 16: Cash from Operating Activities
    # Old: =SUM(col6:col16) -> circular reference (includes its own row)
    # New: =col6+col8+SUM(col10:col15)
    $ws.Range("$col`16").Formula = "=$col`6+$col`8+SUM($col`10:$col`15)"

    # Row 27: Net Change in Cash
    # Old: =col17+col21+col26 (blank rows)
    # New: =col16+col20+col25 (actual subtotal rows)
    $ws.Range("$col`27").Formula = "=$col`16+$col`20+$col`25"
}
